# Regenerate orders with updated distance/size codes.
# Mapping applied to the textual "code" fragments embedded in several columns:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 / S25 / NULL / Face codes are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Columns (1-based) that contain text needing the substitution:
#  2  -> Condition
#  4  -> Filename_Left
#  5  -> Filename_Right
#  8  -> Distance
# 10  -> Size
$textCols = 2, 4, 5, 8, 10

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $textCols) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        if (($old -ne $null) -and ($old -is [string])) {
            $new = $old.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($new -ne $old) {
                $cell.Value = $new
            }
        }
    }
}
